# GAM-305 Traceability Matrix update: add "Level" / "Gameplay/Meta" rows
# (EndGoal win/lose condition logic) to the Traceability Matrix sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Traceability Matrix")

# --- Row 27 : Level / Player start point -------------------------------
# Donor for A=s18, B:H=s19 -> existing row 6 (A6:H6)
$ws.Range("A6:H6").Copy()
$ws.Range("A27:H27").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A27").Value = "Level"
$ws.Range("B27").Value = "Player start point"
$ws.Range("C27").Value = "location for the start of the level"
$ws.Range("D27").Value = "Unreal Engine 5.7.1"
$ws.Range("E27").Value = "alpha"
$ws.Range("F27").Value = "alpha"
$ws.Range("G27").Value = "Bri"
$ws.Range("H27").Value = "Y"

# --- Row 28 : Player exit/end point -------------------------------------
# Donor for B:I=s19 -> existing row 7 (B7:I7)
$ws.Range("B7:I7").Copy()
$ws.Range("B28:I28").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B28").Value = "Player exit/end point"
$ws.Range("C28").Value = "location where player ends/finishes the level"
$ws.Range("D28").Value = "Unreal Engine 5.7.1"
$ws.Range("E28").Value = "alpha"
$ws.Range("F28").Value = "alpha"
$ws.Range("G28").Value = "Bri"
$ws.Range("H28").Value = "Y"
$ws.Range("I28").Value = "Added win/lose condition logic for kittenCollectibles held by the player. If kittens >= 10 then WIN, if kittens < 10 then LOSE. No UI as of yet. Simple print strings for programming testing purposes."

# --- Row 30 : Gameplay/Meta / win condition -----------------------------
# Donor for A=s18, B:I=s19 -> existing row 6 (A6:I6)
$ws.Range("A6:I6").Copy()
$ws.Range("A30:I30").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A30").Value = "Gameplay/Meta"
$ws.Range("B30").Value = "win condition"
$ws.Range("C30").Value = "player must obtain 10 kitens (collectible pickups) throughout the level before reaching the end goal and without getting caught by enemies or side-tracked by obstacles "
$ws.Range("D30").Value = "Unreal Engine 5.7.1"
$ws.Range("E30").Value = "alpha"
$ws.Range("F30").Value = "alpha"
$ws.Range("G30").Value = "Bri"
$ws.Range("H30").Value = "Y"
$ws.Range("I30").Value = "Print string in place for verfication. UI to be added later."

# --- Row 31 : lose condition (kittens < 10) ------------------------------
# Donor for B:I=s19 -> existing row 7 (B7:I7)
$ws.Range("B7:I7").Copy()
$ws.Range("B31:I31").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B31").Value = "lose condition (kittens < 10)"
$ws.Range("C31").Value = "player makes it to end point without all the kittens, gets caught by enemies, and gets side-tracked by obstacle"
$ws.Range("D31").Value = "Unreal Engine 5.7.1"
$ws.Range("E31").Value = "alpha"
$ws.Range("F31").Value = "alpha"
$ws.Range("G31").Value = "Bri"
$ws.Range("H31").Value = "Y"
$ws.Range("I31").Value = "Print string in place for verfication. UI to be added later."

# --- Row 32 : lose condition (caught by customer) ------------------------
# Donor for B:H=s20 -> existing row 17 (B17:H17)
$ws.Rows.Item(32).ClearFormats()
$ws.Range("B17:H17").Copy()
$ws.Range("B32:H32").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B32").Value = "lose condition (caught by customer)"
$ws.Range("C32").Value = "When player is caught by customer enemy, player loses, game is paused, and given option to restart level."
$ws.Range("D32").Value = "Unreal Engine 5.7.1"
$ws.Range("E32").Value = "beta"
$ws.Range("F32").Value = "beta"
$ws.Range("G32").Value = "Bri"
$ws.Range("H32").Value = "N"

# --- Row 33 : lose condition (disciplined by café worker) ----------------
# Donor for B:H=s20 -> existing row 17 (B17:H17)
$ws.Rows.Item(33).ClearFormats()
$ws.Range("B17:H17").Copy()
$ws.Range("B33:H33").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B33").Value = "lose condition (disciplined by café worker)"
$ws.Range("C33").Value = "When player is caught by café worker enemy, player loses, game is paused, and given option to restart level."
$ws.Range("D33").Value = "Unreal Engine 5.7.1"
$ws.Range("E33").Value = "beta"
$ws.Range("F33").Value = "beta"
$ws.Range("G33").Value = "Bri"
$ws.Range("H33").Value = "N"

$excel.CutCopyMode = $false

# --- Selection left where the author last clicked -----------------------
$ws.Range("F37").Select()
